# The whole document body lives in a single <w:p> that uses manual line
# breaks (<w:br/>) instead of paragraph marks. Two spots in that paragraph
# have a "double" line break (an empty <w:r/> run immediately followed by
# two back-to-back <w:r><w:br/></w:r> runs) right before the
# "See plans for location..." sentence - once after the first
# "Utilities Req'd:" line, and once after the second one.
#
# The target edit turns each of those two spots into a real paragraph
# break: the FIRST of the two <w:br/> runs is removed and a paragraph mark
# is inserted in its place, while the SECOND <w:br/> run is left untouched
# (it becomes the first run of the newly created paragraph). This splits
# the single giant paragraph into three paragraphs.

$d = $word.ActiveDocument

function Split-AtUtilitiesDoubleBreak($searchFrom) {
    $full = $d.Content.Text

    # Locate the next "Utilities Req'd:" label ...
    $labelPos = $full.IndexOf("Utilities Req'd:", $searchFrom)

    # ... then walk forward to the first manual line break (chr(11) / "^l")
    # that follows it - that is the first run of the <w:br/><w:br/> pair.
    $breakPos = $labelPos
    while ([int][char]$full[$breakPos] -ne 11) {
        $breakPos = $breakPos + 1
    }

    # Insert a paragraph mark right after that first break run ...
    $breakRange = $d.Range($breakPos, $breakPos + 1)
    $breakRange.InsertParagraphAfter()

    # ... then delete the original break character itself (it now sits
    # just before the freshly inserted paragraph mark), leaving the second
    # <w:br/> run as-is to start the new paragraph.
    $deleteRange = $d.Range($breakPos, $breakPos + 1)
    $deleteRange.Text = ""

    return $labelPos
}

$firstLabelPos = Split-AtUtilitiesDoubleBreak(0)
Split-AtUtilitiesDoubleBreak($firstLabelPos + 1) | Out-Null
